# "added moves and backup database"
# Adds the Ice-type moves to the "Moves" sheet and updates the active
# sheet/selection so "Moves" becomes the visible/active tab.

$wb  = $excel.ActiveWorkbook
$wsPokemon = $wb.Worksheets.Item("Pokemon")
$wsMoves   = $wb.Worksheets.Item("Moves")

# --- Append the new Ice-type move rows (32-43) -----------------------
# Column B ("Ice") is written before column A on the first new row so the
# shared-string table registers "Ice" ahead of the move names, and the
# move-name cells are written in the same order the original author must
# have used (32,33,34,36,35,37,38,39,40,41,42,43) so the shared-string
# indices line up with the target workbook.

$wsMoves.Range("B32").Value = "Ice"
$wsMoves.Range("A32").Value = "Aurora Beam"
$wsMoves.Range("A33").Value = "Blizzard"
$wsMoves.Range("A34").Value = "Hail"
$wsMoves.Range("A36").Value = "Ice Ball"
$wsMoves.Range("A35").Value = "Haze"
$wsMoves.Range("A37").Value = "Ice Beam"
$wsMoves.Range("A38").Value = "Ice Punch"
$wsMoves.Range("A39").Value = "Icicle Spear"
$wsMoves.Range("A40").Value = "Icy Wind"
$wsMoves.Range("A41").Value = "Mist"
$wsMoves.Range("A42").Value = "Powder Snow"
$wsMoves.Range("A43").Value = "Sheer Cold"

$rows = @(
    @{ Row = 32; B = "Ice"; C = 20; D = 65;     E = 100; F = "S" },
    @{ Row = 33; B = "Ice"; C = 10; D = 120;    E = 70;  F = "S" },
    @{ Row = 34; B = "Ice"; C = 10; D = 0;      E = 100; F = "O" },
    @{ Row = 35; B = "Ice"; C = 30; D = 0;      E = 100; F = "O" },
    @{ Row = 36; B = "Ice"; C = 20; D = 30;     E = 90;  F = "P" },
    @{ Row = 37; B = "Ice"; C = 10; D = 95;     E = 100; F = "S" },
    @{ Row = 38; B = "Ice"; C = 15; D = 75;     E = 100; F = "P" },
    @{ Row = 39; B = "Ice"; C = 30; D = 10;     E = 100; F = "P" },
    @{ Row = 40; B = "Ice"; C = 15; D = 55;     E = 95;  F = "S" },
    @{ Row = 41; B = "Ice"; C = 30; D = 0;      E = 100; F = "O" },
    @{ Row = 42; B = "Ice"; C = 25; D = 40;     E = 100; F = "S" },
    @{ Row = 43; B = "Ice"; C = 5;  D = "inf";  E = 30;  F = "S" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $wsMoves.Range("B$n").Value = $r.B
    $wsMoves.Range("C$n").Value = $r.C
    $wsMoves.Range("D$n").Value = $r.D
    $wsMoves.Range("E$n").Value = $r.E
    $wsMoves.Range("F$n").Value = $r.F
}

# --- Update view state: Moves becomes the active/selected sheet ------
$wsPokemon.Range("G96").Select()
$wsMoves.Range("U28").Select()
$wsMoves.Activate()

Write-Output "done"
